# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 845
    4  = 2160
    5  = 57
    6  = 12520
    7  = 12520
    13 = 935
    14 = 13640
    15 = 13932
    17 = 165
    20 = 1047
    23 = 449
    24 = 5030
    25 = 246
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
